$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CONFIG")

# The DSP "Seven 4 Logistics" email list (row 6, column B) drops the
# crystalleparent@gmail.com address -- the cell text is shortened, but the
# existing mailto hyperlink keeps showing the old, full address as its
# display text (Excel leaves the hyperlink's screen text alone when only
# the underlying cell text is edited).
$ws.Range("B6").Value = "mcgarrybecca02@gmail.com; obrien.gene100@gmail.com"

foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$6') {
        $hl.TextToDisplay = "mcgarrybecca02@gmail.com; obrien.gene100@gmail.com;crystalleparent@gmail.com"
    }
}

# Update the saved view state: scrolled so column B is left-most and the
# active selection sits on B22.
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("B22").Select()
